$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 34, pushing the existing rows 34-100 down to 36-102.
$ws.Rows("34:35").Insert()

# Row 34: new "Provincia del Elquí" / Red Globe / Primera entry (10-kilo tray) dated 44580.
$ws.Cells.Item(34, 1).Value = 11
$ws.Cells.Item(34, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(34, 3).Value = "Bíobío"
$ws.Cells.Item(34, 4).Value = 44580
$ws.Cells.Item(34, 5).Value = 8
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100109
$ws.Cells.Item(34, 8).Value = "Uva"
$ws.Cells.Item(34, 9).Value = 100109001
$ws.Cells.Item(34, 10).Value = "Uva"
$ws.Cells.Item(34, 11).Value = "Red Globe"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 150
$ws.Cells.Item(34, 14).Value = 10000
$ws.Cells.Item(34, 15).Value = 11000
$ws.Cells.Item(34, 16).Value = 10467
$ws.Cells.Item(34, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(34, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(34, 19).Value = 1047
$ws.Cells.Item(34, 20).Value = 10

# Row 35: new "Provincia del Elquí" / Superior Seedless / Primera entry (10-kilo tray) dated 44580.
$ws.Cells.Item(35, 1).Value = 11
$ws.Cells.Item(35, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(35, 3).Value = "Bíobío"
$ws.Cells.Item(35, 4).Value = 44580
$ws.Cells.Item(35, 5).Value = 8
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100109
$ws.Cells.Item(35, 8).Value = "Uva"
$ws.Cells.Item(35, 9).Value = 100109001
$ws.Cells.Item(35, 10).Value = "Uva"
$ws.Cells.Item(35, 11).Value = "Superior Seedless"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 220
$ws.Cells.Item(35, 14).Value = 8000
$ws.Cells.Item(35, 15).Value = 8500
$ws.Cells.Item(35, 16).Value = 8227
$ws.Cells.Item(35, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(35, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(35, 19).Value = 823
$ws.Cells.Item(35, 20).Value = 10
